$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-13): date serial value changes from 45175 to 45183
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 45183
}
